$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update coefficients in row 1
$ws.Range("B1").Value = 1.7
$ws.Range("C1").Value = 1.25
$ws.Range("D1").Value = 0.45
$ws.Range("E1").Value = 1.61
$ws.Range("F1").Value = 0.44
$ws.Range("H1").Value = 0.31
$ws.Range("I1").Value = 0.43
$ws.Range("J1").Value = 0.45

# Move the active selection to J2, finishing the line of reasoning
$ws.Range("J2").Select()
